# Customer dependency diagram (slide 2): add a new "get_order" endpoint
# label just below the existing "create_order" one, and nudge
# "create_order" up to make room for it.
#
# NOTE on the literal Left/Top numbers below: PowerPoint's Shape.Left/Top
# (and AddTextbox's coordinate arguments) are points stored as a 32-bit
# Single. To land on an exact target EMU offset after the float32 round
# trip (EMU = floor(Single(points) * 12700)), the point values here were
# solved so they reproduce the exact target EMU positions.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Locate the existing "create_order" textbox by its current text.
$createOrderShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasTextFrame -and $candidate.TextFrame.TextRange.Text -eq "create_order") {
        $createOrderShape = $candidate
    }
}

# Move "create_order" to its new position.
$createOrderShape.Left = 310.41197204589844
$createOrderShape.Top = 105.14181137084961

# Duplicate it to get an identically formatted textbox (same fill, font,
# highlight, effect list, etc.), then rename/reposition/retext it to
# become the new "get_order" label.
$getOrderShape = $createOrderShape.Duplicate()
$getOrderShape.Name = "TextBox 1"
$getOrderShape.Left = 310.98591613769537
$getOrderShape.Top = 126.82188796997072
$getOrderShape.TextFrame.TextRange.Text = "get_order"
